# Restore full translation-input UI: the "Words" and "Progress" sheets were
# missing rows for the word "ਗੁਬਾਰੀ" (mis-indented early return had been
# skipping the write-back). Re-add them here.

$wb = $excel.ActiveWorkbook

# --- Sheet "Words": add row 3 for the word "ਗੁਬਾਰੀ" ---
$wsWords = $wb.Worksheets.Item("Words")

$wsWords.Cells.Item(3, 1).Value = "ਗੁਬਾਰੀ"
$wsWords.Cells.Item(3, 2).Value = "ਗੁਬਾਰੀ"
$wsWords.Cells.Item(3, 3).Value = $true
$wsWords.Cells.Item(3, 4).Value = 45914.27813299769
$wsWords.Cells.Item(3, 5).Value = $true
$wsWords.Cells.Item(3, 6).Value = 45914.27813299769
$wsWords.Cells.Item(3, 7).Value = $true
$wsWords.Cells.Item(3, 8).Value = 45914.27817406424
$wsWords.Cells.Item(3, 9).Value = $false
$wsWords.Cells.Item(3, 10).Value = ""
$wsWords.Cells.Item(3, 11).Value = 0
$wsWords.Cells.Item(3, 12).Value = ""

$wsWords.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWords.Range("F3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWords.Range("H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Sheet "Progress": add rows 8-10 for the word "ਗੁਬਾਰੀ" ---
$wsProgress = $wb.Worksheets.Item("Progress")

$wsProgress.Cells.Item(8, 1).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(8, 2).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(8, 3).Value = ""
$wsProgress.Cells.Item(8, 4).Value = "ਮਨਮੁਖ ਦੁਬਿਧਾ ਦੁਰਮਤਿ ਬਿਆਪੇ ਜਿਨ ਅੰਤਰਿ ਮੋਹ ਗੁਬਾਰੀ ॥"
$wsProgress.Cells.Item(8, 5).Value = 507
$wsProgress.Cells.Item(8, 6).Value = $true
$wsProgress.Cells.Item(8, 7).Value = 45914.27813299769
$wsProgress.Cells.Item(8, 8).Value = "not started"
$wsProgress.Cells.Item(8, 9).Value = ""
$wsProgress.Cells.Item(8, 10).Value = ""
$wsProgress.Cells.Item(8, 11).Value = ""

$wsProgress.Cells.Item(9, 1).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(9, 2).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(9, 3).Value = ""
$wsProgress.Cells.Item(9, 4).Value = "ਜਿਨ੍ਹ੍ਹਿ ਕੀਏ ਤਿਸਹਿ ਨ ਜਾਣਨੀ ਮਨਮੁਖਿ ਗੁਬਾਰੀ ॥"
$wsProgress.Cells.Item(9, 5).Value = 788
$wsProgress.Cells.Item(9, 6).Value = $true
$wsProgress.Cells.Item(9, 7).Value = 45914.27813299769
$wsProgress.Cells.Item(9, 8).Value = "not started"
$wsProgress.Cells.Item(9, 9).Value = ""
$wsProgress.Cells.Item(9, 10).Value = ""
$wsProgress.Cells.Item(9, 11).Value = ""

$wsProgress.Cells.Item(10, 1).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(10, 2).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(10, 3).Value = ""
$wsProgress.Cells.Item(10, 4).Value = "ਬਾਹਰਿ ਭਸਮ ਲੇਪਨ ਕਰੇ ਅੰਤਰਿ ਗੁਬਾਰੀ ॥"
$wsProgress.Cells.Item(10, 5).Value = 1243
$wsProgress.Cells.Item(10, 6).Value = $true
$wsProgress.Cells.Item(10, 7).Value = 45914.27813299769
$wsProgress.Cells.Item(10, 8).Value = "not started"
$wsProgress.Cells.Item(10, 9).Value = ""
$wsProgress.Cells.Item(10, 10).Value = ""
$wsProgress.Cells.Item(10, 11).Value = ""

$wsProgress.Range("G8:G10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
